$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 26.86490966666667
$ws.Range("H2").Value = 80.594729
$ws.Range("I2").Value = 0.1447302967754861
$ws.Range("J2").Value = 0.1447302967754861
$ws.Range("O2").Value = 0.02773017886769741
$ws.Range("P2").Value = 0.02773017886769741
$ws.Range("Q2").Value = 1.395399227966222
$ws.Range("R2").Value = 12.558593051696
$ws.Range("S2").Value = 0.004013397017159161
$ws.Range("T2").Value = 0.004013397017159161

$ws.Range("G3").Value = 26.86490966666667
$ws.Range("H3").Value = 80.594729
$ws.Range("I3").Value = 0.1447302967754861
$ws.Range("J3").Value = 0.1447302967754861
$ws.Range("M3").Value = 1.821156333333333
$ws.Range("N3").Value = 5.463469
$ws.Range("O3").Value = 0.9722698211323025
$ws.Range("P3").Value = 0.9722698211323026
$ws.Range("Q3").Value = 48.92520038387789
$ws.Range("R3").Value = 440.326803454901
$ws.Range("S3").Value = 0.140716899758327
$ws.Range("T3").Value = 0.140716899758327

$ws.Range("I4").Value = 0.4077186109324291
$ws.Range("J4").Value = 0.4077186109324292
$ws.Range("O4").Value = 0.02773017886769741
$ws.Range("P4").Value = 0.02773017886769741
$ws.Range("S4").Value = 0.01130611000884539
$ws.Range("T4").Value = 0.01130611000884539

$ws.Range("I5").Value = 0.4077186109324291
$ws.Range("J5").Value = 0.4077186109324292
$ws.Range("M5").Value = 1.821156333333333
$ws.Range("N5").Value = 5.463469
$ws.Range("O5").Value = 0.9722698211323025
$ws.Range("P5").Value = 0.9722698211323026
$ws.Range("Q5").Value = 137.8268074102651
$ws.Range("R5").Value = 1240.441266692386
$ws.Range("S5").Value = 0.3964125009235837
$ws.Range("T5").Value = 0.3964125009235838

$ws.Range("G6").Value = 14.45399366666666
$ws.Range("H6").Value = 43.36198099999999
$ws.Range("I6").Value = 0.07786852138807973
$ws.Range("J6").Value = 0.07786852138807973
$ws.Range("O6").Value = 0.02773017886769741
$ws.Range("P6").Value = 0.02773017886769741
$ws.Range("Q6").Value = 0.750759703038222
$ws.Range("R6").Value = 6.756837327343998
$ws.Range("S6").Value = 0.002159308026254572
$ws.Range("T6").Value = 0.002159308026254572

$ws.Range("G7").Value = 14.45399366666666
$ws.Range("H7").Value = 43.36198099999999
$ws.Range("I7").Value = 0.07786852138807973
$ws.Range("J7").Value = 0.07786852138807973
$ws.Range("M7").Value = 1.821156333333333
$ws.Range("N7").Value = 5.463469
$ws.Range("O7").Value = 0.9722698211323025
$ws.Range("P7").Value = 0.9722698211323026
$ws.Range("Q7").Value = 26.32298210800989
$ws.Range("R7").Value = 236.906838972089
$ws.Range("S7").Value = 0.07570921336182515
$ws.Range("T7").Value = 0.07570921336182516

$ws.Range("G8").Value = 68.62066266666666
$ws.Range("H8").Value = 205.861988
$ws.Range("I8").Value = 0.369682570904005
$ws.Range("J8").Value = 0.369682570904005
$ws.Range("O8").Value = 0.02773017886769741
$ws.Range("P8").Value = 0.02773017886769741
$ws.Range("Q8").Value = 3.564248713123555
$ws.Range("R8").Value = 32.078238418112
$ws.Range("S8").Value = 0.01025136381543829
$ws.Range("T8").Value = 0.01025136381543829

$ws.Range("G9").Value = 68.62066266666666
$ws.Range("H9").Value = 205.861988
$ws.Range("I9").Value = 0.369682570904005
$ws.Range("J9").Value = 0.369682570904005
$ws.Range("M9").Value = 1.821156333333333
$ws.Range("N9").Value = 5.463469
$ws.Range("O9").Value = 0.9722698211323025
$ws.Range("P9").Value = 0.9722698211323026
$ws.Range("Q9").Value = 124.9689544129302
$ws.Range("R9").Value = 1124.720589716372
$ws.Range("S9").Value = 0.3594312070885667
$ws.Range("T9").Value = 0.3594312070885667
